# NYPD 110th Precinct CompStat weekly report - new crime data collected.
# Updates the report header (volume/number, week-of dates) and the weekly
# crime-complaint statistics table (rows 14-29) with freshly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 29   Number  43" -> "Volume 29   Number  44" ---
$ws.Range("A8").Value = "Volume 29   Number  44"

# --- Header: week-of dates ---
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Row 14: Murder ---
$ws.Range("L14").Value = 100

# --- Row 15: Rape ---
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 34
$ws.Range("L15").Value = 41.666666666666
$ws.Range("M15").Value = 54.545454545454
$ws.Range("N15").Value = 21.428571428571

# --- Row 16: Robbery ---
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 160
$ws.Range("F16").Value = 44
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 76
$ws.Range("I16").Value = 352
$ws.Range("J16").Value = 209
$ws.Range("K16").Value = 68.421052631578
$ws.Range("L16").Value = 51.072961373390
$ws.Range("M16").Value = 11.041009463722
$ws.Range("N16").Value = -73.252279635258

# --- Row 17: Fel. Assault ---
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 30
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 452
$ws.Range("J17").Value = 340
$ws.Range("K17").Value = 32.941176470588
$ws.Range("L17").Value = 48.196721311475
$ws.Range("M17").Value = 117.307692307692
$ws.Range("N17").Value = 7.109004739336

# --- Row 18: Burglary ---
# C18 used to hold the text "0" - it now becomes a real numeric count,
# so restore the normal numeric style (matching the other rows) as well.
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -87.5
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 134
$ws.Range("J18").Value = 145
$ws.Range("K18").Value = -7.586206896551
$ws.Range("L18").Value = -22.543352601156
$ws.Range("M18").Value = -43.933054393305
$ws.Range("N18").Value = -92.864749733759

# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = -51.428571428571
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 107
$ws.Range("H19").Value = -15.887850467289
$ws.Range("I19").Value = 1024
$ws.Range("J19").Value = 543
$ws.Range("K19").Value = 88.581952117863
$ws.Range("L19").Value = 147.941888619855
$ws.Range("M19").Value = 120.689655172414
$ws.Range("N19").Value = -4.029990627928

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -19.047619047619
$ws.Range("I20").Value = 187
$ws.Range("J20").Value = 159
$ws.Range("K20").Value = 17.610062893081
$ws.Range("L20").Value = 48.412698412698
$ws.Range("M20").Value = 46.09375
$ws.Range("N20").Value = -90.345895715023

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 204
$ws.Range("G21").Value = 205
$ws.Range("H21").Value = -0.487804878048
$ws.Range("I21").Value = 2189
$ws.Range("J21").Value = 1433
$ws.Range("K21").Value = 52.756454989532
$ws.Range("L21").Value = 71.417384494909
$ws.Range("M21").Value = 58.508327299058
$ws.Range("N21").Value = -67.191247002398

# --- Row 22: Transit ---
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 33.333333333333
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 80
$ws.Range("I22").Value = 41
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = 70.833333333333
$ws.Range("L22").Value = 241.666666666667
$ws.Range("M22").Value = 64

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = 27.906976744186
$ws.Range("F24").Value = 188
$ws.Range("G24").Value = 174
$ws.Range("H24").Value = 8.045977011494
$ws.Range("I24").Value = 2012
$ws.Range("J24").Value = 1527
$ws.Range("K24").Value = 31.761624099541
$ws.Range("L24").Value = 79.162956366874
$ws.Range("M24").Value = 32.368421052631

# --- Row 25: Misd. Assault ---
$ws.Range("C25").Value = 18
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = 20.634920634920
$ws.Range("I25").Value = 773
$ws.Range("J25").Value = 662
$ws.Range("K25").Value = 16.767371601208
$ws.Range("L25").Value = 29.697986577181
$ws.Range("M25").Value = 40.801457194899

# --- Row 26: UCR Rape* ---
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 7
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 75
$ws.Range("I26").Value = 52
$ws.Range("J26").Value = 52
$ws.Range("L26").Value = 30

# --- Row 27: Other Sex Crimes ---
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 400
$ws.Range("F27").Value = 18
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 112
$ws.Range("J27").Value = 85
$ws.Range("K27").Value = 31.764705882352
$ws.Range("L27").Value = 60

# --- Row 28: Shooting Vic. ---
$ws.Range("L28").Value = -23.529411764705
$ws.Range("N28").Value = -76.363636363636

# --- Row 29: Shooting Inc. ---
$ws.Range("L29").Value = -33.333333333333
$ws.Range("N29").Value = -83.333333333333
